# Update "want to go" counts (col F) and a couple of cover-image URLs
# (col I) that changed between two scrapes of the same event listing.
# These edits show up on both the "展览" (sheet1) and "全部类型" (sheet4)
# sheets, which mirror the same rows with a one-row offset.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (row -> new F value) ---
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    6  = 26
    7  = 574
    8  = 45
    9  = 8336
    10 = 777
    11 = 304
    12 = 1119
    13 = 869
    14 = 64
    15 = 40
    16 = 216
    17 = 136
    18 = 57
    19 = 215
    20 = 908
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}
$ws1.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"

# --- Sheet "全部类型" (row -> new F value) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    7  = 26
    9  = 574
    10 = 45
    11 = 8336
    12 = 777
    13 = 304
    14 = 1119
    15 = 869
    16 = 64
    17 = 40
    18 = 216
    19 = 136
    20 = 57
    21 = 215
    22 = 909
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"
